# 小组计划实施表 — add "第十一周 周三" (week 11, Wednesday) section
# and fill in the completion ("完成情况") column for the previous
# "第十一周 周一&周二" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Fill in the "完成情况" (completion status) column for the
#    existing 第十一周 周一&周二 block (rows 204-209).
# ---------------------------------------------------------------
$ws.Range("C204").Value = "完成“商品类别管理”模块"
$ws.Range("C205").Value = "未完成"
$ws.Range("C206").Value = "完成"
$ws.Range("C207").Value = "未完成"
$ws.Range("C208").Value = "未完成"
$ws.Range("C209").Value = "完成"

# Row 204's wrapped text now needs a touch more height.
$ws.Rows.Item(204).RowHeight = 68

# Update that week's summary line.
$ws.Range("A210").Value = "总结：前端可以参考别人源码"

# ---------------------------------------------------------------
# 2. Append the new 第十一周 周三 block (rows 212-220), reusing the
#    same layout/formatting as the other weekly blocks.
# ---------------------------------------------------------------
$template = $ws.Range("A192:D200")
$template.Copy()
$ws.Range("A212").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A212").Value = "日期：2017.11.8 第十一周 周三"

$ws.Range("A213").Value = "人员"
$ws.Range("B213").Value = "计划任务"
$ws.Range("C213").Value = "完成情况"
$ws.Range("D213").Value = "备注"

$ws.Range("A214").Value = "李杰"
$ws.Range("B214").Value = "开发web app接口所有模块"

$ws.Range("A215").Value = "周振朋"
$ws.Range("B215").Value = "完善首页所有模块，并尝试开发“买卖”模块"

$ws.Range("A216").Value = "禤锦辉"
$ws.Range("B216").Value = "帮助前端开发人员开发其中一个小模块"

$ws.Range("A217").Value = "柯新钿"
$ws.Range("B217").Value = "开发前端“我的”模块中的“我的收藏”模块，并按照“我的收藏”模板做好“我的订单”，“我的发布”，“我的求购”等三个模块"

$ws.Range("A218").Value = "冯文雄"
$ws.Range("B218").Value = "开发web app接口所有模块"

$ws.Range("A219").Value = "阿卜力孜"
$ws.Range("B219").Value = "帮助前端开发人员开发其中一个小模块"

$ws.Range("A220").Value = "总结："

# Row heights for the new block (text wraps to different numbers of
# lines depending on content length).
$ws.Rows.Item(212).RowHeight = 22.5
$ws.Rows.Item(213).RowHeight = 22.5
$ws.Rows.Item(214).RowHeight = 22.5
$ws.Rows.Item(215).RowHeight = 45
$ws.Rows.Item(216).RowHeight = 45
$ws.Rows.Item(217).RowHeight = 112.5
$ws.Rows.Item(218).RowHeight = 22.5
$ws.Rows.Item(219).RowHeight = 45
$ws.Rows.Item(220).RowHeight = 22.5

# Merge the title/banner rows, matching the other weekly blocks.
$ws.Range("A212:D212").Merge()
$ws.Range("A220:D220").Merge()

# ---------------------------------------------------------------
# 3. Misc view state (selection / scroll / window size) to mirror
#    what a human editing this in Excel would have left behind.
# ---------------------------------------------------------------
$ws.Range("C215").Select()
try { $excel.ActiveWindow.ScrollRow = 209 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A209") } catch {}
try {
    $excel.ActiveWindow.Width = 20385
    $excel.ActiveWindow.Height = 8370
} catch {}
